# Remove the "4. Method: Possible algorithms" slide (the table comparing
# algorithm results) from the presentation.
#
# In the original deck this is slide 16 (1-indexed), right before the
# closing "Questions and comments?" slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$s.Delete()
